# Edit: add "2022-Q3" sheet (feat: add 2022-Q3 data)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q3" and shift
#    the existing "2022-Q2"/"2022-Q1" rows down, fixing up the A-column
#    sequence index as we go.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

# Give the new row 2's A cell the same style as A3 (bold/bordered index
# style) without minting a new style entry.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 13
$summary.Cells.Item(2,4).Value = 1.03

# Row that used to be "2022-Q2" (was row 2) is now row 3 - fix its index.
$summary.Cells.Item(3,1).Value = 1

# Row that used to be "2022-Q1" (was row 3) is now row 4 - fix its index.
$summary.Cells.Item(4,1).Value = 2

# ---------------------------------------------------------------------
# 2) New "2022-Q3" sheet: duplicate the "2022-Q2" sheet's layout (so the
#    header/index-column styling matches) immediately before it, rename
#    it, then overwrite the data with the Q3 fund-holdings table.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$qdata = @"
0	000242	景顺长城策略精选	11.62	90.68	5.88	0.6833	2
1	010779	西部利得量化优选一年持有期混合A	4.88	88.29	1.85	0.0903	8
2	202019	南方策略优化混合	2.88	93.97	2.27	0.0654	7
3	001703	银华沪港深增长股票A	2.07	88.44	3.14	0.0650	8
4	005914	景顺长城智能生活混合	1.14	85.26	3.30	0.0376	9
5	006225	人保量化基本面混合A	0.47	90.98	5.32	0.0250	1
6	010780	西部利得量化优选一年持有期混合C	1.22	88.29	1.85	0.0226	8
7	007903	长城量化小盘股票	0.86	91.70	1.24	0.0107	5
8	001744	诺安进取回报灵活配置混合	0.23	82.31	4.30	0.0099	3
9	011231	光大保德信锦弘混合A	1.95	26.05	0.49	0.0096	7
10	014364	银华沪港深增长股票C	0.25	88.44	3.14	0.0078	8
11	011232	光大保德信锦弘混合C	0.82	26.05	0.49	0.0040	7
12	006226	人保量化基本面混合C	0.04	90.98	5.32	0.0021	1
"@

$rows = $qdata -split "`r?`n"

# Stamp the index-column (A) style down to every data row first (copied
# from the existing A2 cell) so every row gets the bordered/bold index
# style the source sheet uses, then fill in the real values below.
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H14").PasteSpecial(-4122)

$r = 2
foreach ($line in $rows) {
    if ($line.Trim().Length -eq 0) { continue }
    $f = $line -split "`t"
    $q3.Cells.Item($r,1).Value = [double]$f[0]
    $q3.Cells.Item($r,2).Value = "'" + $f[1]
    $q3.Cells.Item($r,3).Value = $f[2]
    $q3.Cells.Item($r,4).Value = "'" + $f[3]
    $q3.Cells.Item($r,5).Value = "'" + $f[4]
    $q3.Cells.Item($r,6).Value = "'" + $f[5]
    $q3.Cells.Item($r,7).Value = "'" + $f[6]
    $q3.Cells.Item($r,8).Value = [double]$f[7]
    $r = $r + 1
}
